$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price / volume(1h) data (scraper refresh)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-5.84%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.79%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.044"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.66%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07914"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.89%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.947"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-9.21%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.729"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.91%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.021"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.89%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.873"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.61%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9233"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.43%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1200"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "18.67%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1842"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.46%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09372"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.02%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03535"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.63%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09874"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.53%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001389"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.03%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005884"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.69%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.14%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3444"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.11%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1308"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.25%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.036"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.36%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04488"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.37%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004574"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.31%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001249"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.92%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.92%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01901"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.43%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04703"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.57%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007588"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.20%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009551"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "22.32%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.51%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002109"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.61%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01121"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-7.57%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006006"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-7.44%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.06%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.06%"
